# Fruta / hortaliza, semanal
#
# Insert one new weekly price-report row for "Terminal Hortofrutícola Agro
# Chillán" (Naranja, Valencia, Primera) above the existing row 299, pushing
# the former rows 299-346 down to 300-347.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 299 - this shifts rows 299:346 down to 300:347
# (carrying their values/formatting with them, same as a manual
# right-click > Insert in Excel).
$ws.Rows.Item(299).Insert()

# Populate the newly inserted row 299 with the new weekly record.
$ws.Range("A299").Value = 7
$ws.Range("B299").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C299").Value = "Ñuble"
$ws.Range("D299").Value = 44637
$ws.Range("E299").Value = 16
$ws.Range("F299").Value = "Fruta"
$ws.Range("G299").Value = 100102
$ws.Range("H299").Value = "Cítricos"
$ws.Range("I299").Value = 100102005
$ws.Range("J299").Value = "Naranja"
$ws.Range("K299").Value = "Valencia"
$ws.Range("L299").Value = "Primera"
$ws.Range("M299").Value = 120
$ws.Range("N299").Value = 9000
$ws.Range("O299").Value = 10000
$ws.Range("P299").Value = 9500
$ws.Range("Q299").Value = "$/bandeja 15 kilos granel"
$ws.Range("R299").Value = "Región de O'Higgins"
$ws.Range("S299").Value = 633
$ws.Range("T299").Value = 15
